{"js": "// Remove the early draft \"Fire behavior modeling\" subsection (heading +\n// the HIGRAD/FIRETEC paragraph + the blank spacer paragraph + the long\n// \"Using this modeling framework...\" paragraph) that duplicated the later\n// \"Fire behavior modeling\" section under Methods. This restores the\n// document to the version before that section had been drafted in-line.\n\nconst body = context.document.body;\n\n// There are two \"Fire behavior modeling\" headings in the document; the one\n// we need to remove is immediately followed by a paragraph that begins\n// \"To model fire behavior based on estimated fuels structures in reburns,\n// we used the HIGRAD/FIRETEC system...\". Locate it via search so the\n// script is resilient to exact paragraph indices.\nconst headingHits = body.search(\"Fire behavior modeling\", { matchCase: true });\nheadingHits.load(\"items\");\nawait context.sync();\n\nlet targetHeadingPara = null;\nfor (let i = 0; i < headingHits.items.length; i++) {\n  const headingPara = headingHits.items[i].paragraphs.getFirst();\n  const followingPara = headingPara.getNext();\n  followingPara.load(\"text\");\n  await context.sync();\n\n  if (followingPara.text.indexOf(\"To model fire behavior based on estimated fuels structures in reburns\") !== -1) {\n    targetHeadingPara = headingPara;\n    break;\n  }\n}\n\nif (targetHeadingPara) {\n  // Collect the heading paragraph plus the next three paragraphs:\n  //   1. \"To model fire behavior based on estimated fuels structures...\"\n  //   2. an empty spacer paragraph\n  //   3. \"Using this modeling framework, we modeled predicted fire danger...\"\n  const parasToDelete = [targetHeadingPara];\n  let cursor = targetHeadingPara;\n  for (let i = 0; i < 3; i++) {\n    cursor = cursor.getNext();\n    parasToDelete.push(cursor);\n  }\n\n  // Delete from last to first so earlier references stay valid.\n  for (let i = parasToDelete.length - 1; i >= 0; i--) {\n    parasToDelete[i].delete();\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the early draft \"Fire behavior modeling\" subsection (heading +\n# the HIGRAD/FIRETEC paragraph + the blank spacer paragraph + the long\n# \"Using this modeling framework...\" paragraph) that duplicated the later\n# \"Fire behavior modeling\" section under Methods. This restores the\n# document to the version before that section had been drafted in-line.\n\n$d = $word.ActiveDocument\n\n# There are two \"Fire behavior modeling\" headings in the document; the one\n# we need to remove is immediately followed by a paragraph that begins\n# \"To model fire behavior based on estimated fuels structures in reburns,\n# we used the HIGRAD/FIRETEC system...\". Locate it by walking paragraphs so\n# the script is resilient to exact paragraph indices.\n$targetIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $ptext = $p.Range.Text.Trim()\n    if ($ptext -eq \"Fire behavior modeling\" -and $i -lt $count) {\n        $nextText = $d.Paragraphs.Item($i + 1).Range.Text\n        if ($nextText -like \"*To model fire behavior based on estimated fuels structures in reburns*\") {\n            $targetIndex = $i\n            break\n        }\n    }\n}\n\nif ($targetIndex -gt 0) {\n    # The heading paragraph plus the next three paragraphs:\n    #   1. \"To model fire behavior based on estimated fuels structures...\"\n    #   2. an empty spacer paragraph\n    #   3. \"Using this modeling framework, we modeled predicted fire danger...\"\n    $startPara = $d.Paragraphs.Item($targetIndex)\n    $endPara = $d.Paragraphs.Item($targetIndex + 3)\n\n    $deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)\n    $deleteRange.Delete()\n}\n"}
